# "New Pubs Edit 9-24-19"
#
# SharePoint's Document-ID feature (the "Document Library Form" content-type
# part and the "_dlc_DocId" property-bag part) re-shuffles the two small
# custom XML parts it owns every time the deck is re-published through the
# portal. One part carries the <FormTemplates> display/edit/new form
# pointers, the other carries the <p:properties>/<documentManagement>
# _dlc_DocId / _dlc_DocIdUrl bag. After a republish their package slots
# trade places (and their datastoreItem GUIDs/schemaRefs follow along),
# even though neither part's logical content changes.
#
# Reproduce that swap through CustomXMLParts (the supported automation
# surface for these SharePoint-managed parts) by exchanging the XML payload
# of the two parts in place.

$p = $ppt.ActivePresentation
$parts = $p.CustomXMLParts

$formsNamespace = "http://schemas.microsoft.com/sharepoint/v3/contenttype/forms"
$propsNamespace = "http://schemas.microsoft.com/office/2006/metadata/properties"

$formsXml = "<?mso-contentType ?>" + `
  "<FormTemplates xmlns=`"$formsNamespace`">" + `
  "<Display>DocumentLibraryForm</Display>" + `
  "<Edit>DocumentLibraryForm</Edit>" + `
  "<New>DocumentLibraryForm</New>" + `
  "</FormTemplates>"

$propsXml = "<p:properties xmlns:p=`"$propsNamespace`" " + `
  "xmlns:xsi=`"http://www.w3.org/2001/XMLSchema-instance`" " + `
  "xmlns:pc=`"http://schemas.microsoft.com/office/infopath/2007/PartnerControls`">" + `
  "<documentManagement>" + `
  "<_dlc_DocId xmlns=`"037063e9-a85e-4c78-8627-f1a9315663e5`">EVEA5JW6U4JV-6-9956</_dlc_DocId>" + `
  "<_dlc_DocIdUrl xmlns=`"037063e9-a85e-4c78-8627-f1a9315663e5`">" + `
  "<Url>https://portal.roitraining.com/Courses/_layouts/DocIdRedir.aspx?ID=EVEA5JW6U4JV-6-9956</Url>" + `
  "<Description>EVEA5JW6U4JV-6-9956</Description>" + `
  "</_dlc_DocIdUrl>" + `
  "<Date_x0020_last_x0020_used xmlns=`"027ed24f-5970-4294-be5c-0919c5aaa214`" xsi:nil=`"true`"/>" + `
  "<Customization_x0020_Information xmlns=`"027ed24f-5970-4294-be5c-0919c5aaa214`" xsi:nil=`"true`"/>" + `
  "</documentManagement>" + `
  "</p:properties>"

# Find the two SharePoint-owned parts by the namespace of their root
# element (robust to whatever slot/index they currently occupy) and swap
# their payloads.
$formsPart = $null
$propsPart = $null

for ($i = 1; $i -le $parts.Count; $i++) {
    $part = $parts.Item($i)
    $ns = $null
    try { $ns = $part.NamespaceURI } catch { $ns = $null }
    if ($ns -eq $formsNamespace) {
        $formsPart = $part
    } elseif ($ns -eq $propsNamespace) {
        $propsPart = $part
    }
}

if (($formsPart -ne $null) -and ($propsPart -ne $null)) {
    $formsPart.XML = $propsXml
    $propsPart.XML = $formsXml
} else {
    # Fall back to locating the parts by scanning their current XML text
    # (covers hosts that don't expose NamespaceURI on CustomXMLPart).
    $formsPart = $null
    $propsPart = $null
    for ($i = 1; $i -le $parts.Count; $i++) {
        $part = $parts.Item($i)
        $xml = $part.XML
        if ($xml -ne $null) {
            if ($xml.Contains("FormTemplates")) {
                $formsPart = $part
            } elseif ($xml.Contains("documentManagement")) {
                $propsPart = $part
            }
        }
    }
    if (($formsPart -ne $null) -and ($propsPart -ne $null)) {
        $formsPart.XML = $propsXml
        $propsPart.XML = $formsXml
    }
}
